$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new record row right above the current row 596, pushing the
# existing rows 596:682 down to 597:683 (and extending the used range to
# A1:R683).
$ws.Rows.Item(596).Insert()

# Populate the newly-inserted row with the new price observation.
$ws.Cells.Item(596, 1).Value = 4
$ws.Cells.Item(596, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(596, 3).Value = "Los Lagos"
$ws.Cells.Item(596, 4).Value = 45077
$ws.Cells.Item(596, 5).Value = 10
$ws.Cells.Item(596, 6).Value = 100114001
$ws.Cells.Item(596, 7).Value = "Papa"
$ws.Cells.Item(596, 8).Value = "Red Lady"
$ws.Cells.Item(596, 9).Value = "1a (guarda)"
$ws.Cells.Item(596, 10).Value = 150
$ws.Cells.Item(596, 11).Value = 12000
$ws.Cells.Item(596, 12).Value = 12000
$ws.Cells.Item(596, 13).Value = 12000
$ws.Cells.Item(596, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(596, 15).Value = "Provincia de Llanquihue"
$ws.Cells.Item(596, 16).Value = 480
$ws.Cells.Item(596, 17).Value = 25
$ws.Cells.Item(596, 18).Value = "Hortaliza"
